$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.1351698189973831
$ws.Range("B1").Value = 0.1156130135059357
$ws.Range("C1").Value = 0.1065426841378212
$ws.Range("D1").Value = 0.1216452568769455
$ws.Range("E1").Value = 0.1593261808156967
